$wb = $excel.ActiveWorkbook

# y_fitted_on_begin_2016
$ws = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws.Range("A2").Value = 1995
$ws.Range("B2").Value = 62.7522224925761
$ws.Range("A3").Value = 1996
$ws.Range("B3").Value = 63.62734645199382
$ws.Range("A4").Value = 1997
$ws.Range("B4").Value = 64.15978510197377
$ws.Range("A5").Value = 1998
$ws.Range("B5").Value = 64.89944831178339
$ws.Range("A6").Value = 1999
$ws.Range("B6").Value = 65.91522182477088
$ws.Range("A7").Value = 2000
$ws.Range("B7").Value = 66.4578417426288
$ws.Range("A8").Value = 2001
$ws.Range("B8").Value = 66.7975432524639
$ws.Range("A9").Value = 2002
$ws.Range("B9").Value = 67.64452755037614
$ws.Range("A10").Value = 2003
$ws.Range("B10").Value = 68.56883756341452
$ws.Range("A11").Value = 2004
$ws.Range("B11").Value = 69.57708614403776
$ws.Range("A12").Value = 2005
$ws.Range("B12").Value = 69.41174236615628
$ws.Range("A13").Value = 2006
$ws.Range("B13").Value = 69.54364537854501
$ws.Range("A14").Value = 2007
$ws.Range("B14").Value = 70.02139177040874
$ws.Range("A15").Value = 2008
$ws.Range("B15").Value = 70.43395347937252
$ws.Range("A16").Value = 2009
$ws.Range("B16").Value = 71.02759690222443
$ws.Range("A17").Value = 2010
$ws.Range("B17").Value = 71.38769141891304
$ws.Range("A18").Value = 2011
$ws.Range("B18").Value = 71.50283115190551
$ws.Range("A19").Value = 2012
$ws.Range("B19").Value = 71.31134399330145
$ws.Range("A20").Value = 2013
$ws.Range("B20").Value = 71.19075670774615
$ws.Range("A21").Value = 2014
$ws.Range("B21").Value = 70.6633332140057
$ws.Range("A22").Value = 2015
$ws.Range("B22").Value = 70.23017643232332
$ws.Range("A23").Value = 2016
$ws.Range("B23").Value = 70.1066824489499

# y_pred_on_2017_2021
$ws = $wb.Worksheets.Item("y_pred_on_2017_2021")
$ws.Range("A2").Value = 2017
$ws.Range("B2").Value = 70.14058843822477
$ws.Range("A3").Value = 2018
$ws.Range("B3").Value = 69.91011651629316
$ws.Range("A4").Value = 2019
$ws.Range("B4").Value = 69.71461862479464
$ws.Range("A5").Value = 2020
$ws.Range("B5").Value = 69.48554616072676
$ws.Range("A6").Value = 2021
$ws.Range("B6").Value = 69.23853630483706

# y_fitted_on_begin_2021
$ws = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws.Range("A2").Value = 1995
$ws.Range("B2").Value = 62.81541038917255
$ws.Range("A3").Value = 1996
$ws.Range("B3").Value = 63.26773725397905
$ws.Range("A4").Value = 1997
$ws.Range("B4").Value = 64.12350959645227
$ws.Range("A5").Value = 1998
$ws.Range("B5").Value = 65.15294837797941
$ws.Range("A6").Value = 1999
$ws.Range("B6").Value = 65.8859659106756
$ws.Range("A7").Value = 2000
$ws.Range("B7").Value = 66.55574349269897
$ws.Range("A8").Value = 2001
$ws.Range("B8").Value = 66.99623627802173
$ws.Range("A9").Value = 2002
$ws.Range("B9").Value = 67.85499294058552
$ws.Range("A10").Value = 2003
$ws.Range("B10").Value = 68.71768619215513
$ws.Range("A11").Value = 2004
$ws.Range("B11").Value = 69.18957364246475
$ws.Range("A12").Value = 2005
$ws.Range("B12").Value = 69.36817517943865
$ws.Range("A13").Value = 2006
$ws.Range("B13").Value = 69.7789138979664
$ws.Range("A14").Value = 2007
$ws.Range("B14").Value = 69.8219477659957
$ws.Range("A15").Value = 2008
$ws.Range("B15").Value = 70.22755562574173
$ws.Range("A16").Value = 2009
$ws.Range("B16").Value = 70.67533821186578
$ws.Range("A17").Value = 2010
$ws.Range("B17").Value = 71.55579498650447
$ws.Range("A18").Value = 2011
$ws.Range("B18").Value = 71.53141988794056
$ws.Range("A19").Value = 2012
$ws.Range("B19").Value = 71.42875443995015
$ws.Range("A20").Value = 2013
$ws.Range("B20").Value = 70.81661484965103
$ws.Range("A21").Value = 2014
$ws.Range("B21").Value = 70.75376244744837
$ws.Range("A22").Value = 2015
$ws.Range("B22").Value = 70.68953673738817
$ws.Range("A23").Value = 2016
$ws.Range("B23").Value = 69.94820087794699
$ws.Range("A24").Value = 2017
$ws.Range("B24").Value = 69.93782972968754
$ws.Range("A25").Value = 2018
$ws.Range("B25").Value = 69.67841523695367
$ws.Range("A26").Value = 2019
$ws.Range("B26").Value = 69.48935922232687
$ws.Range("A27").Value = 2020
$ws.Range("B27").Value = 69.14780100093196
$ws.Range("A28").Value = 2021
$ws.Range("B28").Value = 70.68898038798828

# y_pred_on_2022_2026
$ws = $wb.Worksheets.Item("y_pred_on_2022_2026")
$ws.Range("A2").Value = 2022
$ws.Range("B2").Value = 71.6848195330414
$ws.Range("A3").Value = 2023
$ws.Range("B3").Value = 71.89450736423359
$ws.Range("A4").Value = 2024
$ws.Range("B4").Value = 72.04072272716483
$ws.Range("A5").Value = 2025
$ws.Range("B5").Value = 72.29525381878764
$ws.Range("A6").Value = 2026
$ws.Range("B6").Value = 72.62541127535778

# y_fitted_on_begin_2021 shrank from 28 data rows (2-29) to 27 data rows (2-28);
# clear the now-removed trailing row so the sheet dimension shrinks to A1:D28
$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws3.Range("A29:D29").Value = $null

"Edits applied"
